# Apply the daily cryptos price/volume refresh (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.019.73"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.828.42"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.04"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6230"
$ws.Range("E6").Value = "  -5.52%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.90"
$ws.Range("E8").Value = "  +7.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07539"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2907"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.71"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07634"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "1.837.13"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.949"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6635"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.17"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009074"
$ws.Range("E17").Value = "  +5.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.977"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").Value = "28.800.37"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "224.28"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.186"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.72"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.395"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.495"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.021"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.038"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05193"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.841"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7288"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "1.278.32"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.758"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.396"
$ws.Range("E41").Value = "  +5.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8925"
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.45"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Value = "1.979.37"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "63.41"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07410"
$ws.Range("E49").Value = "  -16.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3972"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.863"
$ws.Range("E51").Value = "  +1.44%  "
